$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = "2023-03-06"
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100107
$ws.Range("H18").Value = "Otros"
$ws.Range("I18").Value = 100107011
$ws.Range("J18").Value = "Tuna"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Provincia de Los Andes"
$ws.Range("S18").Value = 1111
$ws.Range("T18").Value = 18
